$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume/1h (E) columns per latest crypto data refresh
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.868.82'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.491.57'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.89'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.06'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.49%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.25%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.097.02'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.07'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.902.27'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.480.45'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.05'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '394.65'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.29'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.23'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.996'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.46%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.35%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.70'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.65%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '162.74'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.70%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.40%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.92'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.98%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0740'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '27.24'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.827.77'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.20'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.80'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.52%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '335.96'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.68'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.18%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.41'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.92%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.97%  '
